$d = $word.ActiveDocument
$d.Content.Find.Execute("6+40=46", $true, $false, $false, $false, $false, $true, 1, $false, "53+37=90", 2) | Out-Null
$d.Content.Find.Execute("33+55=88", $true, $false, $false, $false, $false, $true, 1, $false, "58+16=74", 2) | Out-Null
$d.Content.Find.Execute("22+18=40", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=19", 2) | Out-Null
$d.Content.Find.Execute("60-14=46", $true, $false, $false, $false, $false, $true, 1, $false, "62-34=28", 2) | Out-Null
$d.Content.Find.Execute("16+37=53", $true, $false, $false, $false, $false, $true, 1, $false, "89-74=15", 2) | Out-Null
$d.Content.Find.Execute("34-31=3", $true, $false, $false, $false, $false, $true, 1, $false, "29+24=53", 2) | Out-Null
$d.Content.Find.Execute("80-62=18", $true, $false, $false, $false, $false, $true, 1, $false, "57+2=59", 2) | Out-Null
$d.Content.Find.Execute("50+18=68", $true, $false, $false, $false, $false, $true, 1, $false, "9+12=21", 2) | Out-Null
$d.Content.Find.Execute("50+8=58", $true, $false, $false, $false, $false, $true, 1, $false, "36+13=49", 2) | Out-Null
$d.Content.Find.Execute("31-15=16", $true, $false, $false, $false, $false, $true, 1, $false, "54+36=90", 2) | Out-Null
$d.Content.Find.Execute("87+5=92", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=57", 2) | Out-Null
$d.Content.Find.Execute("35+51=86", $true, $false, $false, $false, $false, $true, 1, $false, "66-45=21", 2) | Out-Null
$d.Content.Find.Execute("44-18=26", $true, $false, $false, $false, $false, $true, 1, $false, "71-16=55", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=24", 2) | Out-Null
$d.Content.Find.Execute("4+93=97", $true, $false, $false, $false, $false, $true, 1, $false, "45+26=71", 2) | Out-Null
$d.Content.Find.Execute("6+61=67", $true, $false, $false, $false, $false, $true, 1, $false, "14+69=83", 2) | Out-Null
$d.Content.Find.Execute("92-22=70", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=21", 2) | Out-Null
$d.Content.Find.Execute("36-10=26", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=9", 2) | Out-Null
$d.Content.Find.Execute("93-67=26", $true, $false, $false, $false, $false, $true, 1, $false, "66-29=37", 2) | Out-Null
$d.Content.Find.Execute("31+65=96", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=60", 2) | Out-Null
$d.Content.Find.Execute("66-6=60", $true, $false, $false, $false, $false, $true, 1, $false, "33+24=57", 2) | Out-Null
$d.Content.Find.Execute("47-33=14", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=83", 2) | Out-Null
$d.Content.Find.Execute("80-14=66", $true, $false, $false, $false, $false, $true, 1, $false, "72-41=31", 2) | Out-Null
$d.Content.Find.Execute("81-35=46", $true, $false, $false, $false, $false, $true, 1, $false, "41+49=90", 2) | Out-Null
$d.Content.Find.Execute("44-42=2", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=22", 2) | Out-Null
$d.Content.Find.Execute("59-30=29", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=95", 2) | Out-Null
$d.Content.Find.Execute("88-5=83", $true, $false, $false, $false, $false, $true, 1, $false, "29+63=92", 2) | Out-Null
$d.Content.Find.Execute("37+19=56", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("82-71=11", $true, $false, $false, $false, $false, $true, 1, $false, "10+3=13", 2) | Out-Null
$d.Content.Find.Execute("19+28=47", $true, $false, $false, $false, $false, $true, 1, $false, "59-5=54", 2) | Out-Null
$d.Content.Find.Execute("68-30=38", $true, $false, $false, $false, $false, $true, 1, $false, "97-15=82", 2) | Out-Null
$d.Content.Find.Execute("27-5=22", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=80", 2) | Out-Null
$d.Content.Find.Execute("49-42=7", $true, $false, $false, $false, $false, $true, 1, $false, "95-83=12", 2) | Out-Null
$d.Content.Find.Execute("58-56=2", $true, $false, $false, $false, $false, $true, 1, $false, "68-50=18", 2) | Out-Null
$d.Content.Find.Execute("47+47=94", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=67", 2) | Out-Null
$d.Content.Find.Execute("47-40=7", $true, $false, $false, $false, $false, $true, 1, $false, "4+13=17", 2) | Out-Null
$d.Content.Find.Execute("40+21=61", $true, $false, $false, $false, $false, $true, 1, $false, "63-53=10", 2) | Out-Null
$d.Content.Find.Execute("7-1=6", $true, $false, $false, $false, $false, $true, 1, $false, "62-53=9", 2) | Out-Null
$d.Content.Find.Execute("10+82=92", $true, $false, $false, $false, $false, $true, 1, $false, "38+21=59", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $false, $false, $false, $false, $true, 1, $false, "72-47=25", 2) | Out-Null
$d.Content.Find.Execute("68-24=44", $true, $false, $false, $false, $false, $true, 1, $false, "4+43=47", 2) | Out-Null
$d.Content.Find.Execute("81-5=76", $true, $false, $false, $false, $false, $true, 1, $false, "73-12=61", 2) | Out-Null
$d.Content.Find.Execute("56+17=73", $true, $false, $false, $false, $false, $true, 1, $false, "56+31=87", 2) | Out-Null
$d.Content.Find.Execute("45-40=5", $true, $false, $false, $false, $false, $true, 1, $false, "49-0=49", 2) | Out-Null
$d.Content.Find.Execute("16+61=77", $true, $false, $false, $false, $false, $true, 1, $false, "44-41=3", 2) | Out-Null
$d.Content.Find.Execute("85+4=89", $true, $false, $false, $false, $false, $true, 1, $false, "73+18=91", 2) | Out-Null
$d.Content.Find.Execute("71-69=2", $true, $false, $false, $false, $false, $true, 1, $false, "40+51=91", 2) | Out-Null
$d.Content.Find.Execute("37+32=69", $true, $false, $false, $false, $false, $true, 1, $false, "89-53=36", 2) | Out-Null
$d.Content.Find.Execute("55-38=17", $true, $false, $false, $false, $false, $true, 1, $false, "34+35=69", 2) | Out-Null
$d.Content.Find.Execute("89-61=28", $true, $false, $false, $false, $false, $true, 1, $false, "20-19=1", 2) | Out-Null
$d.Content.Find.Execute("49-4=45", $true, $false, $false, $false, $false, $true, 1, $false, "40+47=87", 2) | Out-Null
$d.Content.Find.Execute("75-8=67", $true, $false, $false, $false, $false, $true, 1, $false, "71-5=66", 2) | Out-Null
$d.Content.Find.Execute("76-48=28", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("45-24=21", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=48", 2) | Out-Null
$d.Content.Find.Execute("0+24=24", $true, $false, $false, $false, $false, $true, 1, $false, "64-29=35", 2) | Out-Null
$d.Content.Find.Execute("9+20=29", $true, $false, $false, $false, $false, $true, 1, $false, "61-17=44", 2) | Out-Null
$d.Content.Find.Execute("68-9=59", $true, $false, $false, $false, $false, $true, 1, $false, "9+68=77", 2) | Out-Null
$d.Content.Find.Execute("26+20=46", $true, $false, $false, $false, $false, $true, 1, $false, "16+12=28", 2) | Out-Null
$d.Content.Find.Execute("86-17=69", $true, $false, $false, $false, $false, $true, 1, $false, "79+14=93", 2) | Out-Null
$d.Content.Find.Execute("40+37=77", $true, $false, $false, $false, $false, $true, 1, $false, "49-33=16", 2) | Out-Null
$d.Content.Find.Execute("88-24=64", $true, $false, $false, $false, $false, $true, 1, $false, "41+55=96", 2) | Out-Null
$d.Content.Find.Execute("72+2=74", $true, $false, $false, $false, $false, $true, 1, $false, "66+24=90", 2) | Out-Null
$d.Content.Find.Execute("67-24=43", $true, $false, $false, $false, $false, $true, 1, $false, "81-67=14", 2) | Out-Null
$d.Content.Find.Execute("76-13=63", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=9", 2) | Out-Null
$d.Content.Find.Execute("26-23=3", $true, $false, $false, $false, $false, $true, 1, $false, "63-48=15", 2) | Out-Null
$d.Content.Find.Execute("5-3=2", $true, $false, $false, $false, $false, $true, 1, $false, "39+20=59", 2) | Out-Null
$d.Content.Find.Execute("30+36=66", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("32-29=3", $true, $false, $false, $false, $false, $true, 1, $false, "31+60=91", 2) | Out-Null
$d.Content.Find.Execute("94-47=47", $true, $false, $false, $false, $false, $true, 1, $false, "85-22=63", 2) | Out-Null
$d.Content.Find.Execute("85-71=14", $true, $false, $false, $false, $false, $true, 1, $false, "45+29=74", 2) | Out-Null
$d.Content.Find.Execute("9+39=48", $true, $false, $false, $false, $false, $true, 1, $false, "73+2=75", 2) | Out-Null
$d.Content.Find.Execute("84-36=48", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=22", 2) | Out-Null
$d.Content.Find.Execute("28+56=84", $true, $false, $false, $false, $false, $true, 1, $false, "54+36=90", 2) | Out-Null
$d.Content.Find.Execute("63-15=48", $true, $false, $false, $false, $false, $true, 1, $false, "37-31=6", 2) | Out-Null
$d.Content.Find.Execute("14+3=17", $true, $false, $false, $false, $false, $true, 1, $false, "74-26=48", 2) | Out-Null
$d.Content.Find.Execute("2+97=99", $true, $false, $false, $false, $false, $true, 1, $false, "43+47=90", 2) | Out-Null
$d.Content.Find.Execute("40-19=21", $true, $false, $false, $false, $false, $true, 1, $false, "64+33=97", 2) | Out-Null
$d.Content.Find.Execute("65+27=92", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 2) | Out-Null
$d.Content.Find.Execute("31+44=75", $true, $false, $false, $false, $false, $true, 1, $false, "1+76=77", 2) | Out-Null
$d.Content.Find.Execute("9+23=32", $true, $false, $false, $false, $false, $true, 1, $false, "42-14=28", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "0+78=78", 2) | Out-Null
$d.Content.Find.Execute("3+13=16", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=95", 2) | Out-Null
$d.Content.Find.Execute("98-74=24", $true, $false, $false, $false, $false, $true, 1, $false, "53-39=14", 2) | Out-Null
$d.Content.Find.Execute("17+51=68", $true, $false, $false, $false, $false, $true, 1, $false, "2+28=30", 2) | Out-Null
$d.Content.Find.Execute("82-11=71", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=12", 2) | Out-Null
$d.Content.Find.Execute("78-3=75", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=46", 2) | Out-Null
$d.Content.Find.Execute("28+14=42", $true, $false, $false, $false, $false, $true, 1, $false, "35+59=94", 2) | Out-Null
$d.Content.Find.Execute("87-30=57", $true, $false, $false, $false, $false, $true, 1, $false, "90-35=55", 2) | Out-Null
$d.Content.Find.Execute("38+6=44", $true, $false, $false, $false, $false, $true, 1, $false, "8+41=49", 2) | Out-Null
$d.Content.Find.Execute("46-36=10", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=74", 2) | Out-Null
$d.Content.Find.Execute("21+46=67", $true, $false, $false, $false, $false, $true, 1, $false, "70-26=44", 2) | Out-Null
$d.Content.Find.Execute("0+87=87", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=61", 2) | Out-Null
$d.Content.Find.Execute("30+28=58", $true, $false, $false, $false, $false, $true, 1, $false, "36-33=3", 2) | Out-Null
$d.Content.Find.Execute("97-46=51", $true, $false, $false, $false, $false, $true, 1, $false, "83-76=7", 2) | Out-Null
$d.Content.Find.Execute("84-43=41", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=92", 2) | Out-Null
$d.Content.Find.Execute("12-7=5", $true, $false, $false, $false, $false, $true, 1, $false, "65-63=2", 2) | Out-Null
$d.Content.Find.Execute("21+72=93", $true, $false, $false, $false, $false, $true, 1, $false, "37-11=26", 2) | Out-Null
$d.Content.Find.Execute("76-29=47", $true, $false, $false, $false, $false, $true, 1, $false, "16+70=86", 2) | Out-Null
$d.Content.Find.Execute("86-35=51", $true, $false, $false, $false, $false, $true, 1, $false, "44+36=80", 2) | Out-Null
$d.Content.Find.Execute("73-18=55", $true, $false, $false, $false, $false, $true, 1, $false, "48-28=20", 2) | Out-Null
